$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Delta tiempo" column text on existing rows (dedupe "0 min" -> "0 mins") ---
$ws.Range("F11").Value = "0 mins"
$ws.Range("F12").Value = "0 mins"
$ws.Range("F13").Value = "0 mins"

# --- Copy the formatting of row 13 onto the two new rows (14 and 15) so they get the ---
# --- same number formats / styles as the rest of the data rows (drops the unused    ---
# --- "D/M/YYYY" number format that used to be applied to the formerly-blank rows).  ---
$ws.Range("C13:I13").Copy()
$ws.Range("C14:I14").PasteSpecial(-4122)
$ws.Range("C13:I13").Copy()
$ws.Range("C15:I15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 14 ---
$ws.Range("C14").Value = 45346
$ws.Range("D14").Value = "9:47 AM"
$ws.Range("E14").Value = "9:53 AM"
$ws.Range("F14").Value = "0 mins"
$ws.Range("G14").Value = "6 mins"
$ws.Range("H14").Value = "Anexos"
$ws.Range("I14").Value = "Se agregan los formularios LOGT"

# --- Row 15 ---
$ws.Range("C15").Value = 45346
$ws.Range("D15").Value = "9:54 AM"
$ws.Range("E15").Value = "10:28 AM"
$ws.Range("F15").Value = "0 mins"
$ws.Range("G15").Value = "34 mins"
$ws.Range("H15").Value = "Modelo de casos de uso del sistema"
$ws.Range("I15").Value = "Se modifican un par de casos de usos a los roles"
